# Update Name of Algo
# Applies corrected KNN-imputed values to specific cells in column A and C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C3"  = -13.543
    "A12" = -21.721
    "C14" = -12.586
    "C26" = -12.805
    "A27" = -21.908
    "C31" = -13.661
    "A32" = -21.659
    "C35" = -12.762
    "A36" = -20.178
    "C37" = -13.341
    "A38" = -19.741
    "C45" = -12.477
    "A46" = -21.9
    "C52" = -11.587
    "A54" = -22.15
    "A55" = -22.219
    "A56" = -22.158
    "C57" = -13.829
    "A67" = -21.519
    "A69" = -21.721
    "A72" = -21.481
    "C81" = -13.192
    "A83" = -21.627
    "C83" = -12.789
    "A86" = -22.257
    "A91" = -21.522
    "A93" = -21.421
    "A99" = -20.428
    "C100" = -13.018
    "C102" = -12.87
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
